$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# --- Change 1: move "TextBox 2" (the translation-specs box) up slightly ---
# a:off y goes from 585924/2649391 EMU (Top=208.6135pt) to 585924/2605846 EMU (Top=205.1847pt)
$box = $s.Shapes.Item(3)
$box.Top = 205.18474

# --- Change 2: reword + recolor the first run of "TextBox 9" (press release / link line) ---
$linkBox = $s.Shapes.Item(7)
$tr = $linkBox.TextFrame.TextRange
$run1 = $tr.Runs(1)
$run1.Text = "Public visibility: "
$run1.Font.Color.RGB = 8947848   # &H888888 (RGB 136,136,136)
